# On slide 1, the methodology table has a "4a" row (Classification modeling /
# "Is there a way to quantify these characteristics in a predictive model?")
# whose "Data Visualization Approaches" cell was a placeholder "N/A". Redeem
# it now that the preprocessing/model work is done: a Confusion Matrix is the
# visualization used to assess the classification model.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape that holds the methodology table (a graphic frame).
$tableShape = $null
foreach ($shape in $s.Shapes) {
    if ($shape.HasTable) {
        $tableShape = $shape
        break
    }
}

$tbl = $tableShape.Table
$cell = $tbl.Cell(7, 5)
$cell.Shape.TextFrame.TextRange.Text = "Confusion Matrix"
